$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217. This pushes the existing rows 217-309
# down to 218-310 (so the former row 309 becomes row 310), matching the
# dimension change from A1:R309 to A1:R310.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new record.
$ws.Cells.Item(217, 1).Value = 5
$ws.Cells.Item(217, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(217, 3).Value = "Maule"
$ws.Cells.Item(217, 4).Value = 44636
$ws.Cells.Item(217, 5).Value = 7
$ws.Cells.Item(217, 6).Value = 100112032
$ws.Cells.Item(217, 7).Value = "Zapallo italiano"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 350
$ws.Cells.Item(217, 11).Value = 11000
$ws.Cells.Item(217, 12).Value = 11000
$ws.Cells.Item(217, 13).Value = 11000
$ws.Cells.Item(217, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(217, 15).Value = "Región del Maule"
$ws.Cells.Item(217, 16).Value = 220
$ws.Cells.Item(217, 17).Value = 50
$ws.Cells.Item(217, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# rest of the column (style index 2 in the original workbook).
$ws.Cells.Item(217, 4).NumberFormat = $ws.Cells.Item(218, 4).NumberFormat
